$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "button1..4" columns (E:G) entirely from the header/data rows
# and the old per-row numeric "button threshold" columns (D:G on rows 4-7).
$ws.Range("D1:G2").Clear()
$ws.Range("D4:G7").Clear()

# New header row 1 / row 2 content for column D ("预定义" / "const")
$ws.Range("D1").Value = "预定义"
$ws.Range("D2").Value = "const"

# New data rows: replace numeric thresholds with named constants
$ws.Range("D4").Value = "one"
$ws.Range("D5").Value = "two"
$ws.Range("D6").Value = "three"
$ws.Range("D7").Value = "four"

# Update the active selection to match the authored state
$ws.Range("D3").Select()
